$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: direct assignment is safe.
# Numeric-looking Price values (column D) need NumberFormat forced to
# Text ("@") before assignment, otherwise Excel auto-converts the typed
# string into a real number (e.g. "1.00" -> 1), which would silently
# drop the trailing zero / change the stored type. We restore the
# cell style back to Normal/General afterwards so formatting is
# unaffected, while the underlying value remains the exact text.

$ws.Range('D2').Value = '67.377.35'
$ws.Range('E2').Value = '  -0.53%  '
$ws.Range('D3').Value = '3.521.00'
$ws.Range('E3').Value = '  -0.90%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').NumberFormat = 'General'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '611.83'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.42'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.87%  '
$ws.Range('D7').Value = '3.520.59'
$ws.Range('E7').Value = '  -0.80%  '
$ws.Range('E9').Value = '  -0.90%  '
$ws.Range('E10').Value = '  -0.96%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.05'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.425'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.42%  '
$ws.Range('D14').Value = '4.115.13'
$ws.Range('E14').Value = '  -0.94%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '32.01'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.47%  '
$ws.Range('D16').Value = '3.519.59'
$ws.Range('E16').Value = '  -1.55%  '
$ws.Range('D17').Value = '67.376.90'
$ws.Range('E17').Value = '  -0.55%  '
$ws.Range('E18').Value = '  +0.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.41'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.37%  '
$ws.Range('E20').Value = '  -2.24%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '444.77'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.38'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.98%  '
$ws.Range('E23').Value = '  -2.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '77.37'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.44%  '
$ws.Range('E25').Value = '  +9.20%  '
$ws.Range('D26').Value = '3.660.29'
$ws.Range('E26').Value = '  -1.00%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.37'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.79%  '
$ws.Range('E28').Value = '  +0.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.40'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.57%  '
$ws.Range('E30').Value = '  -2.37%  '
$ws.Range('E31').Value = '  -4.69%  '
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('E33').Value = '  +3.63%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.92'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.18'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.00%  '
$ws.Range('D36').Value = '3.510.98'
$ws.Range('E36').Value = '  -1.20%  '
$ws.Range('E37').Value = '  -3.05%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.03'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.48%  '
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '177.94'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.43%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.18'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.83%  '
$ws.Range('E43').Value = '  -0.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.45'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.80%  '
$ws.Range('E45').Value = '  -0.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '28.74'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.85%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '44.80'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.50%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.63'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.50%  '
$ws.Range('E49').Value = '  +4.68%  '
$ws.Range('E50').Value = '  -1.08%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.00'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.77%  '
